$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03096523827643
$ws.Range("D2").Value = 1.035670106207959
$ws.Range("E2").Value = 1.04016315221547
$ws.Range("F2").Value = 1.050173548980149
$ws.Range("I2").Value = 1.036737571951997
$ws.Range("J2").Value = 1.036103572741627
$ws.Range("K2").Value = 1.038466046285523
$ws.Range("L2").Value = 1.042946280839702
$ws.Range("M2").Value = 1.052928561516346
$ws.Range("N2").Value = 1.016033180115432
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031781540173511
$ws.Range("D3").Value = 1.036287812876352
$ws.Range("E3").Value = 1.040962196227817
$ws.Range("F3").Value = 1.051222092908845
$ws.Range("I3").Value = 1.036932112982919
$ws.Range("J3").Value = 1.036562229554795
$ws.Range("K3").Value = 1.038893705444089
$ws.Range("L3").Value = 1.043555721757979
$ws.Range("M3").Value = 1.053788890011594
$ws.Range("N3").Value = 1.016185515656919
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032310119955753
$ws.Range("D4").Value = 1.036687685725202
$ws.Range("E4").Value = 1.041480027669654
$ws.Range("F4").Value = 1.051901904800973
$ws.Range("I4").Value = 1.037056755599795
$ws.Range("J4").Value = 1.036858724338186
$ws.Range("K4").Value = 1.039169912579465
$ws.Range("L4").Value = 1.043950212478039
$ws.Range("M4").Value = 1.054346309372932
$ws.Range("N4").Value = 1.016283970880135
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032532423869958
$ws.Range("D5").Value = 1.036855832586624
$ws.Range("E5").Value = 1.041697913131681
$ws.Range("F5").Value = 1.052188015577204
$ws.Range("I5").Value = 1.037108858123408
$ws.Range("J5").Value = 1.036983300960664
$ws.Range("K5").Value = 1.039285905118198
$ws.Range("L5").Value = 1.044116089058382
$ws.Range("M5").Value = 1.054580821696364
$ws.Range("N5").Value = 1.016325333268882
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032569754849544
$ws.Range("D6").Value = 1.036884067500394
$ws.Range("E6").Value = 1.041734508104224
$ws.Range("F6").Value = 1.052236073432079
$ws.Range("I6").Value = 1.037117588914797
$ws.Range("J6").Value = 1.037004213799173
$ws.Range("K6").Value = 1.039305373414765
$ws.Range("L6").Value = 1.044143942346541
$ws.Range("M6").Value = 1.05462020746551
$ws.Range("N6").Value = 1.016332276531344
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03231309004379
$ws.Range("D7").Value = 1.036689932355469
$ws.Range("E7").Value = 1.041482938324701
$ws.Range("F7").Value = 1.051905726580044
$ws.Range("I7").Value = 1.03705745296512
$ws.Range("J7").Value = 1.036860389212714
$ws.Range("K7").Value = 1.039171462969865
$ws.Range("L7").Value = 1.043952428802161
$ws.Range("M7").Value = 1.054349442259137
$ws.Range("N7").Value = 1.016284523677659
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031241032188169
$ws.Range("D8").Value = 1.035878825769303
$ws.Range("E8").Value = 1.040433027009715
$ws.Range("F8").Value = 1.05052763224074
$ws.Range("I8").Value = 1.036803573981459
$ws.Range("J8").Value = 1.036258636682795
$ws.Range("K8").Value = 1.038610682023415
$ws.Range("L8").Value = 1.043152214166994
$ws.Range("M8").Value = 1.053219161728522
$ws.Range("N8").Value = 1.016084686413746
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029354889127461
$ws.Range("D9").Value = 1.034450968335421
$ws.Range("E9").Value = 1.038589118943525
$ws.Range("F9").Value = 1.048109534785613
$ws.Range("I9").Value = 1.036346754850507
$ws.Range("J9").Value = 1.035196124123926
$ws.Range("K9").Value = 1.03761860754751
$ws.Range("L9").Value = 1.041743273750457
$ws.Range("M9").Value = 1.051233109213306
$ws.Range("N9").Value = 1.015731676518239
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028099538674374
$ws.Range("D10").Value = 1.033500113939974
$ws.Range("E10").Value = 1.037364088257125
$ws.Range("F10").Value = 1.046504469160707
$ws.Range("I10").Value = 1.036035898165679
$ws.Range("J10").Value = 1.034486404432568
$ws.Range("K10").Value = 1.036954668154839
$ws.Range("L10").Value = 1.040804820174183
$ws.Range("M10").Value = 1.049912951607176
$ws.Range("N10").Value = 1.015495774778218
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027556469436907
$ws.Range("D11").Value = 1.03308865340405
$ws.Range("E11").Value = 1.036834661982568
$ws.Range("F11").Value = 1.045811136322757
$ws.Range("I11").Value = 1.035899806271125
$ws.Range("J11").Value = 1.034178774963379
$ws.Range("K11").Value = 1.036666583467217
$ws.Range("L11").Value = 1.040398674062472
$ws.Range("M11").Value = 1.049342245407989
$ws.Range("N11").Value = 1.015393498474148
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027354826795218
$ws.Range("D12").Value = 1.032935860225649
$ws.Range("E12").Value = 1.036638164317478
$ws.Range("F12").Value = 1.045553854152113
$ws.Range("I12").Value = 1.035849032713899
$ws.Range("J12").Value = 1.034064461280493
$ws.Range("K12").Value = 1.03655948783448
$ws.Range("L12").Value = 1.040247846160852
$ws.Range("M12").Value = 1.049130400895868
$ws.Range("N12").Value = 1.015355489474593
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02739807629791
$ws.Range("D13").Value = 1.032968632972538
$ws.Range("E13").Value = 1.036680306702289
$ws.Range("F13").Value = 1.045609030601007
$ws.Range("I13").Value = 1.035859933892454
$ws.Range("J13").Value = 1.034088984035619
$ws.Range("K13").Value = 1.036582464181547
$ws.Range("L13").Value = 1.040280197751621
$ws.Range("M13").Value = 1.049175835852472
$ws.Range("N13").Value = 1.01536364338861
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027539800004746
$ws.Range("D14").Value = 1.033076022614307
$ws.Range("E14").Value = 1.036818416256802
$ws.Range("F14").Value = 1.045789864138668
$ws.Range("I14").Value = 1.035895613857348
$ws.Range("J14").Value = 1.034169326698901
$ws.Range("K14").Value = 1.036657732700325
$ws.Range("L14").Value = 1.040386205903893
$ws.Range("M14").Value = 1.04932473138669
$ws.Range("N14").Value = 1.015390357024216
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027627130965126
$ws.Range("D15").Value = 1.033142194477136
$ws.Range("E15").Value = 1.036903530674378
$ws.Range("F15").Value = 1.045901315153506
$ws.Range("I15").Value = 1.035917567937485
$ws.Range("J15").Value = 1.034218822349689
$ws.Range("K15").Value = 1.036704096472845
$ws.Range("L15").Value = 1.040451525414567
$ws.Range("M15").Value = 1.049416489587643
$ws.Range("N15").Value = 1.015406813667512
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028135591172469
$ws.Range("D16").Value = 1.033527426960888
$ws.Range("E16").Value = 1.037399246161448
$ws.Range("F16").Value = 1.046550518710306
$ws.Range("I16").Value = 1.036044898828747
$ws.Range("J16").Value = 1.034506814219945
$ws.Range("K16").Value = 1.036973774960325
$ws.Range("L16").Value = 1.04083177929318
$ws.Range("M16").Value = 1.049950847222765
$ws.Range("N16").Value = 1.015502559829637
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028454671274912
$ws.Range("D17").Value = 1.033769145431895
$ws.Range("E17").Value = 1.037710469605574
$ws.Range("F17").Value = 1.046958195385669
$ws.Range("I17").Value = 1.036124372021851
$ws.Range("J17").Value = 1.034687380042246
$ws.Range("K17").Value = 1.037142778800318
$ws.Range("L17").Value = 1.04107035972148
$ws.Range("M17").Value = 1.050286285788301
$ws.Range("N17").Value = 1.015562584533396
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028640834003677
$ws.Range("D18").Value = 1.033910161250253
$ws.Range("E18").Value = 1.037892099294425
$ws.Range("F18").Value = 1.047196147439762
$ws.Range("I18").Value = 1.036170583685208
$ws.Range("J18").Value = 1.034792670453008
$ws.Range("K18").Value = 1.037241298436974
$ws.Range("L18").Value = 1.041209539874836
$ws.Range("M18").Value = 1.050482031142459
$ws.Range("N18").Value = 1.015597583425863
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028704318899752
$ws.Range("D19").Value = 1.033958248278624
$ws.Range("E19").Value = 1.037954046953559
$ws.Range("F19").Value = 1.047277310230454
$ws.Range("I19").Value = 1.036186316271276
$ws.Range("J19").Value = 1.034828566542067
$ws.Range("K19").Value = 1.037274881307961
$ws.Range("L19").Value = 1.041257000103792
$ws.Range("M19").Value = 1.05054879042589
$ws.Range("N19").Value = 1.015609515011105
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028420431951461
$ws.Range("D20").Value = 1.033743208675466
$ws.Range("E20").Value = 1.037677068094446
$ws.Range("F20").Value = 1.046914438858031
$ws.Range("I20").Value = 1.036115860159284
$ws.Range("J20").Value = 1.034668010193304
$ws.Range("K20").Value = 1.037124652220948
$ws.Range("L20").Value = 1.041044760198498
$ws.Range("M20").Value = 1.050250287115504
$ws.Range("N20").Value = 1.015556145740979
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027498063755601
$ws.Range("D21").Value = 1.033044397880098
$ws.Range("E21").Value = 1.03677774214874
$ws.Range("F21").Value = 1.045736606186706
$ws.Range("I21").Value = 1.035885113146194
$ws.Range("J21").Value = 1.034145669059122
$ws.Range("K21").Value = 1.036635570421607
$ws.Range("L21").Value = 1.040354988228621
$ws.Range("M21").Value = 1.049280881454881
$ws.Range("N21").Value = 1.015382491046372
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026918583108886
$ws.Range("D22").Value = 1.032605269217951
$ws.Range("E22").Value = 1.036213196784242
$ws.Range("F22").Value = 1.044997517734924
$ws.Range("I22").Value = 1.035738743674165
$ws.Range("J22").Value = 1.033816984539217
$ws.Range("K22").Value = 1.036327555837214
$ws.Range("L22").Value = 1.039921491869583
$ws.Range("M22").Value = 1.048672195138308
$ws.Range("N22").Value = 1.015273197563861
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027225733729108
$ws.Range("D23").Value = 1.032838036190476
$ws.Range("E23").Value = 1.036512387457153
$ws.Range("F23").Value = 1.045389183472696
$ws.Range("I23").Value = 1.035816458904841
$ws.Range("J23").Value = 1.033991251450812
$ws.Range("K23").Value = 1.036490888094834
$ws.Range("L23").Value = 1.040151278010078
$ws.Range("M23").Value = 1.048994793331986
$ws.Range("N23").Value = 1.015331146389549
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028435903077589
$ws.Range("D24").Value = 1.033754928300602
$ws.Range("E24").Value = 1.037692160495965
$ws.Range("F24").Value = 1.046934210051621
$ws.Range("I24").Value = 1.036119706747693
$ws.Range("J24").Value = 1.034676762688882
$ws.Range("K24").Value = 1.037132843019543
$ws.Range("L24").Value = 1.041056327458656
$ws.Range("M24").Value = 1.050266553089134
$ws.Range("N24").Value = 1.015559055192792
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029842142477711
$ws.Range("D25").Value = 1.034819925740068
$ws.Range("E25").Value = 1.039065072626984
$ws.Range("F25").Value = 1.048733443648346
$ws.Range("I25").Value = 1.036465969080741
$ws.Range("J25").Value = 1.03547105680045
$ws.Range("K25").Value = 1.037875538307351
$ws.Range("L25").Value = 1.042107375946529
$ws.Range("M25").Value = 1.051745873933094
$ws.Range("N25").Value = 1.015823038806952
